# Update "Datos actualizados" timestamp (row 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value2 = "Datos actualizados a 28 de Septiembre de 2020 a las 01:13"

# --- Venezuela overtakes Costa Rica (row 54 / row 55 swap position) ---
# Row 54 now holds Venezuela's refreshed numbers
$ws.Range("A54").Value2 = "Venezuela"
$ws.Range("B54").Value2 = 72691
$ws.Range("C54").Value2 = 751
$ws.Range("D54").Value2 = 62427
$ws.Range("E54").Value2 = 9658
$ws.Range("F54").Value2 = 0
$ws.Range("G54").Value2 = 6
$ws.Range("H54").Value2 = 606

# Row 55 now holds Costa Rica's (unchanged/carried-over) numbers
$ws.Range("A55").Value2 = "Costa Rica"
$ws.Range("B55").Value2 = 72049
$ws.Range("C55").Value2 = 0
$ws.Range("D55").Value2 = 27760
$ws.Range("E55").Value2 = 43461
$ws.Range("F55").Value2 = 0
$ws.Range("G55").Value2 = 0
$ws.Range("H55").Value2 = 828

# --- Santa Lucia / Timor Oriental swap position (tied values) ---
$ws.Range("A207").Value2 = "Santa Lucia"
$ws.Range("A208").Value2 = "Timor Oriental"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value2 = 7319656
$ws.Range("C4").Value2 = 32095
$ws.Range("D4").Value2 = 4545936
$ws.Range("E4").Value2 = 2564274
$ws.Range("G4").Value2 = 269
$ws.Range("H4").Value2 = 209446

# --- Brasil (row 6) ---
$ws.Range("D6").Value2 = 4060088
$ws.Range("E6").Value2 = 530480

# --- Colombia (row 8) ---
$ws.Range("B8").Value2 = 813056
$ws.Range("C8").Value2 = 7018
$ws.Range("D8").Value2 = 711472
$ws.Range("E8").Value2 = 76096
$ws.Range("G8").Value2 = 192
$ws.Range("H8").Value2 = 25488

# --- Argentina (row 12) ---
$ws.Range("B12").Value2 = 711325
$ws.Range("C12").Value2 = 8841
$ws.Range("D12").Value2 = 565935
$ws.Range("E12").Value2 = 129641
$ws.Range("G12").Value2 = 206
$ws.Range("H12").Value2 = 15749

# --- Canada (row 29) ---
$ws.Range("B29").Value2 = 153125
$ws.Range("C29").Value2 = 1454
$ws.Range("D29").Value2 = 131098
$ws.Range("E29").Value2 = 12759

# --- Panama (row 38) ---
$ws.Range("B38").Value2 = 110555
$ws.Range("C38").Value2 = 447
$ws.Range("D38").Value2 = 87215
$ws.Range("E38").Value2 = 21000
$ws.Range("G38").Value2 = 17
$ws.Range("H38").Value2 = 2340

# --- Egipto (row 41) ---
$ws.Range("B41").Value2 = 102840
$ws.Range("C41").Value2 = 104
$ws.Range("D41").Value2 = 95080
$ws.Range("E41").Value2 = 1877
$ws.Range("G41").Value2 = 14
$ws.Range("H41").Value2 = 5883

# --- Japon (row 48) ---
$ws.Range("B48").Value2 = 81690
$ws.Range("C48").Value2 = 635
$ws.Range("D48").Value2 = 74607
$ws.Range("E48").Value2 = 5538
$ws.Range("G48").Value2 = 5
$ws.Range("H48").Value2 = 1545

# --- Nigeria (row 58) ---
$ws.Range("B58").Value2 = 58324
$ws.Range("C58").Value2 = 126
$ws.Range("D58").Value2 = 49794
$ws.Range("E58").Value2 = 7422
$ws.Range("G58").Value2 = 2
$ws.Range("H58").Value2 = 1108

# --- Zambia (row 92) ---
$ws.Range("B92").Value2 = 14641
$ws.Range("C92").Value2 = 29
$ws.Range("D92").Value2 = 13784
$ws.Range("E92").Value2 = 525

# --- Consejo Danes para los Refugiados (row 99) ---
$ws.Range("B99").Value2 = 10612
$ws.Range("C99").Value2 = 19
$ws.Range("E99").Value2 = 248

# --- Guinea (row 100) ---
$ws.Range("B100").Value2 = 10580
$ws.Range("C100").Value2 = 68
$ws.Range("D100").Value2 = 9892
$ws.Range("E100").Value2 = 622

# --- Zimbabue (row 112) ---
$ws.Range("B112").Value2 = 7812
$ws.Range("C112").Value2 = 9
$ws.Range("D112").Value2 = 6106
$ws.Range("E112").Value2 = 1479

# --- Mauritania (row 114) ---
$ws.Range("B114").Value2 = 7464
$ws.Range("C114").Value2 = 2
$ws.Range("E114").Value2 = 233

# --- Trinidad yTobago (row 133) ---
$ws.Range("B133").Value2 = 4362
$ws.Range("C133").Value2 = 50
$ws.Range("E133").Value2 = 2057
$ws.Range("G133").Value2 = 1
$ws.Range("H133").Value2 = 71

# --- Uruguay (row 154) ---
$ws.Range("B154").Value2 = 2008
$ws.Range("C154").Value2 = 10
$ws.Range("D154").Value2 = 1728
$ws.Range("E154").Value2 = 233

Write-Host "Edits applied"
